# "gets column headers properly from Tesla and Panera public statements"
#
# The labels in column A of rows 2-28 were off by one row (each caption was
# sitting one row above where it belonged). This restores the correct
# alignment by shifting every caption down one row and inserting the
# missing " net income" caption at the top (row 2). Rows 19-29 also pick up
# their Tesla cash-flow-statement figures (columns B/C/D) which previously
# had no values, row 30 becomes the "noncontrolling interests in
# subsidiaries" balance-sheet line (Panera figures), and row 31 becomes the
# "total liabilities and equity" line - which, unlike the row it displaced,
# only reports two years of data (B/C), so its old third-year (D) figure is
# cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = " net income"

$ws.Range("A3").Value = " depreciation and amortization"

$ws.Range("A4").Value = " stock-based compensation expense"

$ws.Range("A5").Value = " tax benefit from stock-based compensation"

$ws.Range("A6").Value = " other"

$ws.Range("A7").Value = " inventories"

$ws.Range("A8").Value = " prepaid expenses and other"

$ws.Range("A9").Value = " deposits and other"

$ws.Range("A10").Value = " accounts payable"

$ws.Range("A11").Value = " accrued expenses"

$ws.Range("A12").Value = " deferred rent"

$ws.Range("A13").Value = " other net long-term cash provided liabilities by operating activities"

$ws.Range("A14").Value = " additions to property and equipment"

$ws.Range("A15").Value = " proceeds net from cash sale-leaseback used in investing transactions activities"

$ws.Range("A16").Value = " exercise of employee stock options"

$ws.Range("A17").Value = " tax benefit from stock-based compensation"

$ws.Range("A18").Value = " cash and cash equivalents at end of the period accompanying notes are an integral part of the consolidated financial statements."

$ws.Range("A19").Value = " net cash provided by operating activities"
$ws.Range("B19").Value = 5943
$ws.Range("C19").Value = 2405
$ws.Range("D19").Value = 2098

$ws.Range("A20").Value = " purchases of solar energy systems net of sales"
$ws.Range("B20").Value = -75
$ws.Range("C20").Value = -105
$ws.Range("D20").Value = -218

$ws.Range("A21").Value = " business combinations net of cash acquired"
$ws.Range("B21").Value = -13
$ws.Range("C21").Value = -45
$ws.Range("D21").Value = -18

$ws.Range("A22").Value = " net cash used in investing activities"
$ws.Range("B22").Value = -3132
$ws.Range("C22").Value = -1436
$ws.Range("D22").Value = -2337

$ws.Range("A23").Value = " proceeds from issuances of convertible and other debt"
$ws.Range("B23").Value = 9713
$ws.Range("C23").Value = 10669
$ws.Range("D23").Value = 6176

$ws.Range("A24").Value = " repayments of convertible and other debt"
$ws.Range("B24").Value = -11623
$ws.Range("C24").Value = -9161
$ws.Range("D24").Value = -5247

$ws.Range("A25").Value = " collateralized lease repayments"
$ws.Range("B25").Value = -240
$ws.Range("C25").Value = -389
$ws.Range("D25").Value = -559

$ws.Range("A26").Value = " principal payments on finance leases"
$ws.Range("B26").Value = -338
$ws.Range("C26").Value = -321
$ws.Range("D26").Value = -181

$ws.Range("A27").Value = " debt issuance costs"
$ws.Range("B27").Value = -6
$ws.Range("C27").Value = -37
$ws.Range("D27").Value = -15

$ws.Range("A28").Value = " distributions paid to noncontrolling interests in subsidiaries gl)"
$ws.Range("B28").Value = -208
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = -227

$ws.Range("A29").Value = " net cash provided by financing activities"
$ws.Range("B29").Value = 9973
$ws.Range("C29").Value = 1529
$ws.Range("D29").Value = 574

$ws.Range("A30").Value = " noncontrolling interests in subsidiaries"
$ws.Range("B30").Value = 850
$ws.Range("C30").Value = 849

$ws.Range("A31").Value = " total liabilities and equity "
$ws.Range("B31").Value = 52148
$ws.Range("C31").Value = 34309
# row 31 no longer carries a 3rd-year (D) figure, unlike the row it displaced
$ws.Range("D31").ClearContents()

# row 32 (" payments for buy-outs of noncontrolling interests in subsidiaries") is untouched by this edit
